$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block currently occupying A5:L8 needs to move up to A1:L4.
# Capture the existing values (Value2 avoids locale/currency formatting
# issues and works reliably for plain numbers), clear the old range,
# then write the values into the new location.
$srcRange = $ws.Range("A5:L8")
$values = $srcRange.Value2

$srcRange.ClearContents()

$destRange = $ws.Range("A1:L4")
$destRange.Value2 = $values

# Update the active selection to match the new layout.
$ws.Range("L9").Select()
